$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = "Add Position Categories"
$ws.Range("B3").Value = "PASSED"

$ws.Range("A4").Value = "Edit Position Categories"
$ws.Range("B4").Value = "PASSED"

$ws.Range("A5").Value = "Delete Position Categories"
$ws.Range("B5").Value = "PASSED"

$ws.Range("A6").Value = "Add New School Department"
$ws.Range("B6").Value = "PASSED"

$ws.Range("A7").Value = "Edit School Department"
$ws.Range("B7").Value = "PASSED"

$ws.Range("A8").Value = "Delete The School Department"
$ws.Range("B8").Value = "PASSED"
